# ---------------------------------------------------------------------------
# Applies two changes captured by the target commit:
#
# 1. Slide 16's table (the only table in the deck) switches its table style
#    from the custom "Table_0" style {6854F383-44C3-4EA9-95B6-95F076CE4EC9}
#    to the built-in style {536597D9-30D3-4888-BD77-FC01340828D9}.
#
# 2. The presentation's theme colour palette is swapped from the "Integral"
#    palette to the default Office "Office Theme" palette (the commit swaps
#    the raw theme1.xml / theme2.xml package parts; the net, user-visible
#    effect on the live (slide-master-bound) theme is that every themed
#    colour slot changes from the Integral values to the stock Office
#    values).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Retarget the table's style on slide 16 -----------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{536597D9-30D3-4888-BD77-FC01340828D9}")
    }
}

# --- 2. Swap the theme colour scheme back to the stock Office palette ------
# Theme colour slot order for ThemeColorScheme / ColorScheme.Colors:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1 6 accent2 7 accent3 8 accent4
#   9 accent5 10 accent6 11 hlink 12 folHlink
$officeTheme = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeTheme[$i - 1]
}
